$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 13 and 14 (the second and third "dct:creator" rows: Minka and Diba),
# shifting all subsequent rows up by two.
$ws.Range("A13:T14").EntireRow.Delete()
